$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain text (avoids Excel
# auto-converting numeric-looking strings like "581.24" or "1.00" into numbers),
# then restores the original cell style so no formatting is changed.
function Set-TextValue([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "63.218.71"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "3.465.72"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "581.24"
$ws.Range("E5").Value = "  +0.43%  "
Set-TextValue "D6" "147.90"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").Value = "3.465.52"
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("E8").Value = "  -0.06%  "
Set-TextValue "D9" "0.478"
$ws.Range("E9").Value = "  +0.79%  "
Set-TextValue "D10" "7.80"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("E12").Value = "  +5.31%  "
$ws.Range("D13").Value = "4.062.01"
$ws.Range("E13").Value = "  +1.46%  "
Set-TextValue "D14" "29.45"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").Value = "3.459.07"
$ws.Range("E16").Value = "  +1.01%  "
Set-TextValue "D17" "0.0000172"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "63.232.81"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("E19").Value = "  +3.91%  "
Set-TextValue "D20" "14.54"
$ws.Range("E20").Value = "  +3.75%  "
Set-TextValue "D21" "9.32"
$ws.Range("E21").Value = "  +1.81%  "
Set-TextValue "D22" "388.59"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("E23").Value = "  +2.26%  "
Set-TextValue "D24" "74.50"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "3.611.25"
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("E29").Value = "  +2.73%  "
Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  +0.05%  "
Set-TextValue "D31" "8.20"
$ws.Range("E31").Value = "  +2.44%  "
Set-TextValue "D32" "2.14"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D34" "23.46"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D35" "1.34"
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("E37").Value = "  +2.82%  "
Set-TextValue "D38" "1.60"
$ws.Range("E38").Value = "  +5.30%  "
Set-TextValue "D39" "31.89"
$ws.Range("E39").Value = "  +9.64%  "
Set-TextValue "D40" "168.22"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").Value = "3.504.01"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("E42").Value = "  +3.12%  "
Set-TextValue "D43" "0.794"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("E44").Value = "  +4.54%  "
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("E46").Value = "  +3.56%  "
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").Value = "2.591.73"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +10.81%  "
$ws.Range("E50").Value = "  +2.95%  "
Set-TextValue "D51" "23.04"
$ws.Range("E51").Value = "  +0.86%  "
